# "Generate Report for Archive"
# Refresh the localization status report: the handoff status for the
# tracked file has moved from "Ready for handoff" to "In Translation" on
# every sheet that surfaces it (Overview summary columns + each per-locale
# status table), then let the status columns auto-size to the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns("E:F").AutoFit()

# --- zh-cn sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns("C:C").AutoFit()

# --- de-de sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns("C:C").AutoFit()
